$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Columns A and C should default to Text format ("@") so codes with
#    leading zeros (e.g. "007", "012", "01", "03") are preserved as typed.
#    Apply this BEFORE typing the text values so Excel keeps the literal
#    digits instead of parsing them as numbers.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).NumberFormat = "@"
$ws.Range("C1").NumberFormat = "@"
$ws.Range("C3:C4").NumberFormat = "@"

# ---------------------------------------------------------------------------
# 2) Header row (row 1)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "AREA_CODE"
$ws.Range("B1").Value = "AREA_NAME"
$ws.Range("C1").Value = "BRANCH_CODE"
$ws.Range("D1").Value = "BRANCH_NAME"
$ws.Range("E1").Value = "RM_CODE"
$ws.Range("F1").Value = "RM_NAME"
$ws.Range("G1").Value = "BST_CODE"
$ws.Range("H1").Value = "BST_NAME"
$ws.Range("I1").Value = "OS_TARGET_AMT"
$ws.Range("J1").Value = "DISB_TARGET_AMT"
$ws.Range("K1").Value = "INC_TARGET_AMT"

# ---------------------------------------------------------------------------
# 3) Data rows
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "01"
$ws.Range("B2").Value = "Dhaka Area"
$ws.Range("C2").Value = 124
$ws.Range("C2").NumberFormat = "@"
$ws.Range("D2").Value = "Ashkona Branch"
$ws.Range("E2").Value = "RB0506"
$ws.Range("F2").Value = "Kutubuddin Ahmed"
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = 100
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = 50

$ws.Range("A3").Value = "03"
$ws.Range("B3").Value = "Chattogram Area 1"
$ws.Range("C3").Value = "007"
$ws.Range("D3").Value = "Chandgaon Branch"
$ws.Range("E3").Value = "RB0385"
$ws.Range("F3").Value = "Mahabub Hossain"
$ws.Range("G3").Value = ""
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = 50
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = 10

$ws.Range("A4").Value = "03"
$ws.Range("B4").Value = "Chattogram Area 1"
$ws.Range("C4").Value = "012"
$ws.Range("D4").Value = "Feni SME Branch"
$ws.Range("E4").Value = "RB0451"
$ws.Range("F4").Value = "Md. Salahuddin Ahmed"
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = 200
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = 100

# ---------------------------------------------------------------------------
# 4) Formatting clean-up: several data cells (D:J columns in rows 2-4) keep
#    a bold 8pt font from the old template - reset them to the plain body
#    style by copying the already-correct format from column B.
# ---------------------------------------------------------------------------
$ws.Range("B2").Copy() | Out-Null
$ws.Range("D2:J2").PasteSpecial(-4122) | Out-Null
$ws.Range("B3").Copy() | Out-Null
$ws.Range("D3:J3").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Copy() | Out-Null
$ws.Range("D4:J4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 5) Column widths / layout
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 21.7109375
$ws.Columns.Item(6).ColumnWidth = 27.28515625
$ws.Columns.Item(11).ColumnWidth = 21.140625

# ---------------------------------------------------------------------------
# 6) Sheet view: clear frozen/scrolled top-left cell, update selection
# ---------------------------------------------------------------------------
$ws.Range("K14").Select()
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
